$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column D and E to text format to preserve formatting like "42.691.74" and "  +0.81%  "
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "42.691.74"
$ws.Range("E2").Value = "  +0.81%  "
$ws.Range("D3").Value = "2.316.05"
$ws.Range("E3").Value = "  +0.05%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "319.46"
$ws.Range("E5").Value = "  +2.82%  "
$ws.Range("D6").Value = "103.58"
$ws.Range("E6").Value = "  -2.35%  "
$ws.Range("D7").Value = "0.633"
$ws.Range("E7").Value = "  +0.62%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("E9").Value = "  +0.91%  "
$ws.Range("D10").Value = "39.81"
$ws.Range("E10").Value = "  -0.94%  "
$ws.Range("D11").Value = "0.0911"
$ws.Range("E11").Value = "  -0.56%  "
$ws.Range("D12").Value = "8.34"
$ws.Range("E12").Value = "  -0.16%  "
$ws.Range("D13").Value = "0.106"
$ws.Range("E13").Value = "  +0.43%  "
$ws.Range("D14").Value = "0.972"
$ws.Range("E14").Value = "  -0.24%  "
$ws.Range("D15").Value = "15.39"
$ws.Range("E15").Value = "  -1.29%  "
$ws.Range("D16").Value = "2.666.79"
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("D17").Value = "2.315.77"
$ws.Range("E17").Value = "  +0.73%  "
$ws.Range("D18").Value = "42.696.66"
$ws.Range("E18").Value = "  +0.92%  "
$ws.Range("D19").Value = "7.48"
$ws.Range("E19").Value = "  -0.09%  "
$ws.Range("E20").Value = "  +0.92%  "
$ws.Range("E21").Value = "  +4.55%  "
$ws.Range("D22").Value = "73.47"
$ws.Range("E22").Value = "  -2.62%  "
$ws.Range("D23").Value = "279.39"
$ws.Range("E23").Value = "  +6.94%  "
$ws.Range("D24").Value = "10.89"
$ws.Range("E24").Value = "  +17.14%  "
$ws.Range("E25").Value = "  +0.16%  "
$ws.Range("E26").Value = "  -0.40%  "
$ws.Range("D27").Value = "10.91"
$ws.Range("E27").Value = "  -1.41%  "
$ws.Range("D28").Value = "2.36"
$ws.Range("E28").Value = "  +4.31%  "
$ws.Range("D29").Value = "22.98"
$ws.Range("E29").Value = "  -0.47%  "
$ws.Range("D30").Value = "36.01"
$ws.Range("E30").Value = "  +0.42%  "
$ws.Range("D31").Value = "165.09"
$ws.Range("E31").Value = "  -0.39%  "
$ws.Range("D32").Value = "0.0878"
$ws.Range("E32").Value = "  -1.90%  "
$ws.Range("D33").Value = "5.94"
$ws.Range("E33").Value = "  +0.41%  "
$ws.Range("D34").Value = "0.136"
$ws.Range("E34").Value = "  +5.43%  "
$ws.Range("D35").Value = "2.61"
$ws.Range("E35").Value = "  -9.81%  "
$ws.Range("E36").Value = "  -2.06%  "
$ws.Range("D37").Value = "4.62"
$ws.Range("E37").Value = "  +1.78%  "
$ws.Range("D38").Value = "0.0361"
$ws.Range("E38").Value = "  +2.69%  "
$ws.Range("D39").Value = "3.71"
$ws.Range("E39").Value = "  +0.23%  "
$ws.Range("E40").Value = "  +5.20%  "
$ws.Range("E41").Value = "  +2.35%  "
$ws.Range("D42").Value = "99.08"
$ws.Range("E42").Value = "  -0.81%  "
$ws.Range("D43").Value = "69.74"
$ws.Range("E43").Value = "  -0.92%  "
$ws.Range("E44").Value = "  -1.89%  "
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("D46").Value = "12.12"
$ws.Range("E46").Value = "  -0.52%  "
$ws.Range("D47").Value = "113.79"
$ws.Range("E47").Value = "  +1.96%  "
$ws.Range("D48").Value = "80.41"
$ws.Range("E48").Value = "  +8.76%  "
$ws.Range("D49").Value = "9.00"
$ws.Range("E50").Value = "  -1.45%  "
$ws.Range("D51").Value = "1.616.23"
$ws.Range("E51").Value = "  +4.99%  "
